$wb = $excel.ActiveWorkbook

# Data for TODOS / COMBINADAS sheets (identical arrival board)
$rows = @(
    @("15:50", "16_SANTA ANA", 1, "🚌"),
    @("15:53", "11_ETCHEVERRY", 4, "🚌"),
    @("15:56", "27_EL RETIRO", 7, "🚌"),
    @("15:56", "17_ROMERO", 7, "📅"),
    @("16:01", "10_OLMOS", 12, "🚌"),
    @("16:02", "16_SANTA ANA", 13, "🚌"),
    @("16:04", "23_HERNANDEZ", 15, "🚌"),
    @("16:08", "14_ABASTO", 19, "🚌"),
    @("16:13", "215C_LA PLATA", 24, "🚌"),
    @("16:15", "225_C ROCA-H SUR", 26, "📅"),
    @("16:20", "215C_EL PATO", 31, "🚌"),
    @("16:21", "26_HERNANDEZ", 32, "🚌"),
    @("16:29", "10_OLMOS", 40, "🚌"),
    @("16:30", "15_ABASTO", 41, "🚌"),
    @("16:36", "11_ETCHEVERRY", 47, "🚌"),
    @("16:42", "16_P MOR-SANTA ANA", 53, "🚌"),
    @("16:43", "225_GOMEZ", 54, "📅"),
    @("16:48", "15_ABASTO", 59, "🚌"),
    @("16:50", "14_ABASTO", 61, "🚌"),
    @("16:53", "215B_LP-P MOR-40 Y 115", 64, "🚌"),
    @("16:56", "17_179 Y 38", 67, "🚌"),
    @("17:04", "215A_EL PATO", 75, "🚌"),
    @("17:06", "23_HERNANDEZ", 77, "🚌"),
    @("17:14", "215A_LA PLATA", 85, "🚌"),
    @("17:21", "26_HERNANDEZ", 92, "🚌"),
    @("17:24", "84_COLONIA URQUIZA-ESC 49", 95, "🚌"),
    @("17:36", "27_EL RETIRO", 107, "🚌"),
    @("17:38", "17_ROMERO", 109, "🚌"),
    @("17:40", "215B_EL PATO", 111, "📅"),
)

# Data for the "215"-only filtered sheet
$rows215 = @(
    @("16:13", "215C_LA PLATA", 24, "🚌"),
    @("16:20", "215C_EL PATO", 31, "🚌"),
    @("16:53", "215B_LP-P MOR-40 Y 115", 64, "🚌"),
    @("17:04", "215A_EL PATO", 75, "🚌"),
    @("17:14", "215A_LA PLATA", 85, "🚌"),
    @("17:40", "215B_EL PATO", 111, "📅"),
)

function Update-Sheet($ws, $data) {
    for ($i = 0; $i -lt $data.Count; $i++) {
        $r = $i + 2
        $row = $data[$i]
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
    $lastRow = $data.Count + 2
    $ws.Rows.Item($lastRow).Delete()
}

Update-Sheet $wb.Worksheets.Item("TODOS") $rows
Update-Sheet $wb.Worksheets.Item("215") $rows215
Update-Sheet $wb.Worksheets.Item("COMBINADAS") $rows
